$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Update status column (A) for several existing TODO rows.
#    "Not done" -> "Future" / "Done" (reuses existing shared strings).
# ---------------------------------------------------------------------------
$ws.Range("A44").Value = "Future"
$ws.Range("A47").Value = "Future"
$ws.Range("A49").Value = "Future"
$ws.Range("A53").Value = "Done"
$ws.Range("A56").Value = "Future"
$ws.Range("A61").Value = "Done"
$ws.Range("A64").Value = "Future"

# ---------------------------------------------------------------------------
# 2. Apply existing cell formatting (fill/style) to the new rows before
#    filling in their values, so the resulting style indexes match the
#    rest of the "Future" task table (A column centered style, B column
#    highlighted style).
# ---------------------------------------------------------------------------
$ws.Range("A67").Copy()
$ws.Range("A84").PasteSpecial(-4122)

$ws.Range("B67").Copy()
foreach ($r in @(84,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103)) {
    $ws.Range("B$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Populate the new TODO rows (84, 86-103) with values.
#    Column B values are entered in the order they were originally typed so
#    that newly created shared-string entries line up with the source file.
# ---------------------------------------------------------------------------
$ws.Range("B86").Value = "Janine has to fix what Aron broke in the Macros"
$ws.Range("B87").Value = "Janine can sort TODO list"
$ws.Range("B88").Value = "Tower fixes  - Steve & Mike"
$ws.Range("B91").Value = "Molten salt linear fresnel"
$ws.Range("B92").Value = "Reorganize dispatch widget to UI"
$ws.Range("B93").Value = "Sample files - esp Novatec - Mike"
$ws.Range("B94").Value = "Wind - cost and scaling model - Janine"
$ws.Range("B96").Value = "Check the reports"
$ws.Range("B97").Value = "Check the macros"
$ws.Range("B98").Value = "Documentation updates from Paul"
$ws.Range("B99").Value = "Progress updates for solarpilot - Aron"

$ws.Range("C88").Value = "Steve/Mike"
$ws.Range("C89").Value = "Ty/Steve"
$ws.Range("C93").Value = "Mike"

$ws.Range("B89").Value = "Steam tower"
$ws.Range("B90").Value = "Cavity receiver"
$ws.Range("C90").Value = "Ty/Mike/Steve"

$ws.Range("B95").Value = "Review default values, financial, cost #s for PV res/com/util"

$ws.Range("B100").Value = "Subhourly simulation for physical trough"
$ws.Range("C100").Value = "Aron/Mike"

$ws.Range("B101").Value = "Check all results, and summarize for release notes"
$ws.Range("B102").Value = "Add performance adjustment factors to wind model"

$ws.Range("B84").Value = "Curtailment month by hour factors in popup widget thingy"
$ws.Range("B103").Value = "Possible registration issues"

# Remaining column A/C values - these reuse shared strings already present
# in the workbook, so their ordering is not significant.
$ws.Range("A84").Value = "Future"
$ws.Range("C84").Value = "Aron"
$ws.Range("C86").Value = "Janine"
$ws.Range("C87").Value = "Janine"
$ws.Range("C91").Value = "Steve"
$ws.Range("C92").Value = "Steve"
$ws.Range("C94").Value = "Janine"
$ws.Range("C95").Value = "Everyone"
$ws.Range("C96").Value = "Everyone"
$ws.Range("C97").Value = "Everyone"
$ws.Range("C98").Value = "Paul"
$ws.Range("C99").Value = "Aron"
$ws.Range("C101").Value = "Everyone"
$ws.Range("C102").Value = "Janine/Aron"

# ---------------------------------------------------------------------------
# 4. Update the active selection / scroll position to reflect where the
#    sheet was left (bottom of the newly-extended TODO list).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$win.ScrollColumn = 1
$ws.Range("C103").Select()

# ---------------------------------------------------------------------------
# 5. Shrink the workbook window width as recorded in the source workbook.
# ---------------------------------------------------------------------------
$win.Width = 15600
